# Append one row (row 6) of data to Sheet1, as was done at 2025-05-01T15:44:16.289Z
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

# Columns that hold numeric-looking text must be entered with a leading
# apostrophe so Excel stores them as text (not as numbers), matching the
# rest of the sheet where quantities such as "2", "500", "12", "20" are
# all stored as text.
$ws.Cells.Item($row, 1).Value = "'21"
$ws.Cells.Item($row, 2).Value = "احمد"
$ws.Cells.Item($row, 3).Value = "'50"
$ws.Cells.Item($row, 4).Value = "الصمود"
$ws.Cells.Item($row, 5).Value = "الرحلة 1"
$ws.Cells.Item($row, 6).Value = "C1"
$ws.Cells.Item($row, 7).Value = "UNICEF"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٤٤:١٦ م"
